# Update cryptos list price/volume figures to match the latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '43.205.80'
$ws.Range("E2").Value = '  +1.91%  '
$ws.Range("D3").Value = '2.382.48'
$ws.Range("E3").Value = '  +4.11%  '
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("D5").Value = '303.21'
$ws.Range("E5").Value = '  +0.79%  '
$ws.Range("D6").Value = '97.04'
$ws.Range("E6").Value = '  +1.89%  '
$ws.Range("D7").Value = '0.510'
$ws.Range("E7").Value = '  +0.44%  '
$ws.Range("E8").Value = '  -0.10%  '
$ws.Range("E9").Value = '  +2.08%  '
$ws.Range("E10").Value = '  -0.56%  '
$ws.Range("E11").Value = '  +1.34%  '
$ws.Range("E12").Value = '  +2.51%  '
$ws.Range("D13").Value = '18.43'
$ws.Range("E13").Value = '  -2.80%  '
$ws.Range("D14").Value = '6.79'
$ws.Range("E14").Value = '  +1.15%  '
$ws.Range("D15").Value = '2.752.36'
$ws.Range("E15").Value = '  +4.08%  '
$ws.Range("D16").Value = '2.382.03'
$ws.Range("E16").Value = '  +4.35%  '
$ws.Range("D17").Value = '0.809'
$ws.Range("E17").Value = '  +4.12%  '
$ws.Range("D18").Value = '43.184.97'
$ws.Range("E18").Value = '  +1.97%  '
$ws.Range("D19").Value = '12.25'
$ws.Range("E19").Value = '  +0.83%  '
$ws.Range("D20").Value = '6.34'
$ws.Range("E20").Value = '  +6.69%  '
$ws.Range("E21").Value = '  +0.47%  '
$ws.Range("D22").Value = '68.58'
$ws.Range("E22").Value = '  +1.58%  '
$ws.Range("D23").Value = '235.44'
$ws.Range("E23").Value = '  +0.24%  '
$ws.Range("D24").Value = '2.23'
$ws.Range("E24").Value = '  -1.73%  '
$ws.Range("E25").Value = '  +1.80%  '
$ws.Range("E26").Value = '  -0.10%  '
$ws.Range("D27").Value = '24.82'
$ws.Range("E27").Value = '  +2.58%  '
$ws.Range("E28").Value = '  +0.47%  '
$ws.Range("E29").Value = '  +1.23%  '
$ws.Range("D30").Value = '31.59'
$ws.Range("E30").Value = '  -0.40%  '
$ws.Range("E31").Value = '  +0.02%  '
$ws.Range("E32").Value = '  +2.83%  '
$ws.Range("D33").Value = '0.0735'
$ws.Range("E33").Value = '  +6.06%  '
$ws.Range("D34").Value = '17.18'
$ws.Range("E34").Value = '  -1.65%  '
$ws.Range("E35").Value = '  +7.39%  '
$ws.Range("E36").Value = '  +2.76%  '
$ws.Range("E37").Value = '  -0.92%  '
$ws.Range("E38").Value = '  -0.39%  '
$ws.Range("D39").Value = '2.80'
$ws.Range("E39").Value = '  +4.92%  '
$ws.Range("D40").Value = '22.40'
$ws.Range("E40").Value = '  +12.75%  '
$ws.Range("E41").Value = '  +0.57%  '
$ws.Range("D42").Value = '106.21'
$ws.Range("E42").Value = '  -35.63%  '
$ws.Range("D43").Value = '1.956.82'
$ws.Range("E43").Value = '  +0.05%  '
$ws.Range("E44").Value = '  +1.06%  '
$ws.Range("E45").Value = '  +2.69%  '
$ws.Range("E46").Value = '  +0.84%  '
$ws.Range("D47").Value = '9.26'
$ws.Range("E47").Value = '  -10.36%  '
$ws.Range("D48").Value = '52.84'
$ws.Range("E48").Value = '  -0.04%  '
$ws.Range("E49").Value = '  +3.20%  '
$ws.Range("D50").Value = '72.05'
$ws.Range("E50").Value = '  +2.10%  '
$ws.Range("E51").Value = '  +1.61%  '
